$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared-string values, defined here as PowerShell single-quoted here-strings so that
# backticks / quotes inside the Cypher queries are treated completely literally.

$tabStudyFiles = @'
StudyFilesTab
'@

$qSamples = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE diag.stage_of_disease IN ['Unknown']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@

$qStudyFiles = @'
MATCH (f:file)-->(s:study) 
MATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
WHERE diag.stage_of_disease IN ['Unknown']
WITH DISTINCT f,  s, c, demo, diag
WITH
        f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH    
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
WITH    
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$qFiles = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f)-[*]->(samp:sample)
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE diag.stage_of_disease IN ['Unknown'] 
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp:sample)
WITH
        f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN 
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_type, '') AS `File Type`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
'@

$qCases = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
WHERE diag.stage_of_disease IN ['Unknown']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age
RETURN  
       coalesce(c.case_id, '') AS `Case ID`,
       coalesce(s.clinical_study_designation, '') AS `Study Code`,
       coalesce(s.clinical_study_type, '') AS  `Study Type`,
       coalesce(demo.breed, '') AS Breed ,
       coalesce(diag.disease_term, '') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '') AS `Stage of Disease`,
       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
'@

$qStats = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE diag.stage_of_disease IN ['Unknown']  
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# --- Populate the new shared strings in the same order they were introduced in the
#     authored workbook, so the rebuilt shared-string table lines up with the target. ---

# 1) StudyFilesTab label (row 5, col A)
$ws.Range("A5").Value = $tabStudyFiles

# 2) Samples query (row 3, col B)
$ws.Range("B3").Value = $qSamples

# 3) Study files query (row 5, col B)
$ws.Range("B5").Value = $qStudyFiles

# 4) Files query (row 4, col B)
$ws.Range("B4").Value = $qFiles

# 5) Cases query (row 2, col B)
$ws.Range("B2").Value = $qCases

# --- Fill in the rest of the rows (labels / stat query / filenames) ---

$ws.Range("A2").Value = "CasesTab"
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("A4").Value = "FilesTab"

$ws.Range("C2").Value = $qStats
$ws.Range("C3").Value = $qStats
$ws.Range("C4").Value = $qStats
$ws.Range("C5").Value = $qStats

$ws.Range("D2").Value = "TC11_Canine_Filter_StageOfDisease-Unknown_Neo4jData.xlsx"
$ws.Range("D3").Value = "TC11_Canine_Filter_StageOfDisease-Unknown_Neo4jData.xlsx"
$ws.Range("D4").Value = "TC11_Canine_Filter_StageOfDisease-Unknown_Neo4jData.xlsx"
$ws.Range("D5").Value = "TC11_Canine_Filter_StageOfDisease-Unknown_Neo4jData.xlsx"

$ws.Range("E2").Value = "TC11_Canine_Filter_StageOfDisease-Unknown_WebData.xlsx"
$ws.Range("E3").Value = "TC11_Canine_Filter_StageOfDisease-Unknown_WebData.xlsx"
$ws.Range("E4").Value = "TC11_Canine_Filter_StageOfDisease-Unknown_WebData.xlsx"
$ws.Range("E5").Value = "TC11_Canine_Filter_StageOfDisease-Unknown_WebData.xlsx"

# --- Row heights (ht attribute in the sheet XML) ---
# Row 3 keeps its existing height (225) so it is left untouched on purpose.
$ws.Rows.Item(2).RowHeight = 300
$ws.Rows.Item(4).RowHeight = 409.5
$ws.Rows.Item(5).RowHeight = 375

# --- Wrap-text styling on the query columns (B & C), matching the style already used
#     by those columns in the sheet (style index 1 in styles.xml) ---
$ws.Range("B2:C5").WrapText = $true

# --- View state: zoom level and active selection ---
$excel.ActiveWindow.Zoom = 60
[void]$ws.Range("B2").Select()

Write-Host "edit complete"
